# Update cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to remain text,
# matching the original inline-string cell type so numeric-looking
# values like "0.999" or "1.00" are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "63.316.18"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").Value = "2.565.87"
$ws.Range("E3").Value = "  +0.73%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "584.92"
$ws.Range("E5").Value = "  +3.15%  "

$ws.Range("D6").Value = "148.03"
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "0.604"
$ws.Range("E8").Value = "  +4.24%  "

$ws.Range("D10").Value = "5.68"
$ws.Range("E10").Value = "  +1.75%  "

$ws.Range("E11").Value = "  +0.44%  "

$ws.Range("D12").Value = "0.359"
$ws.Range("E12").Value = "  +2.17%  "

$ws.Range("D13").Value = "27.59"
$ws.Range("E13").Value = "  +2.29%  "

$ws.Range("D14").Value = "3.023.36"
$ws.Range("E14").Value = "  +0.66%  "

$ws.Range("D15").Value = "63.174.07"
$ws.Range("E15").Value = "  +0.43%  "

$ws.Range("D16").Value = "0.0000149"
$ws.Range("E16").Value = "  +5.47%  "

$ws.Range("D17").Value = "2.568.69"
$ws.Range("E17").Value = "  +3.68%  "

$ws.Range("D18").Value = "11.39"
$ws.Range("E18").Value = "  -0.71%  "

$ws.Range("D19").Value = "342.98"
$ws.Range("E19").Value = "  +2.73%  "

$ws.Range("E20").Value = "  +4.13%  "

$ws.Range("D21").Value = "6.86"
$ws.Range("E21").Value = "  +1.56%  "

$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").Value = "66.65"
$ws.Range("E23").Value = "  +3.06%  "

$ws.Range("D24").Value = "2.683.78"
$ws.Range("E24").Value = "  +0.15%  "

$ws.Range("E25").Value = "  +3.32%  "

$ws.Range("E26").Value = "  +1.00%  "

$ws.Range("D27").Value = "8.14"
$ws.Range("E27").Value = "  +13.61%  "

$ws.Range("E28").Value = "  +2.80%  "

$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("D31").Value = "1.99"
$ws.Range("E31").Value = "  +8.07%  "

$ws.Range("D32").Value = "0.0₃0830"
$ws.Range("E32").Value = "  +2.60%  "

$ws.Range("D33").Value = "469.26"
$ws.Range("E33").Value = "  +15.55%  "

$ws.Range("D34").Value = "176.72"
$ws.Range("E34").Value = "  -0.18%  "

$ws.Range("E35").Value = "  +3.05%  "

$ws.Range("D36").Value = "0.406"
$ws.Range("E36").Value = "  +2.52%  "

$ws.Range("D37").Value = "19.31"
$ws.Range("E37").Value = "  +2.81%  "

$ws.Range("D38").Value = "4.52"
$ws.Range("E38").Value = "  +4.97%  "

$ws.Range("D40").Value = "1.76"
$ws.Range("E40").Value = "  +0.79%  "

$ws.Range("D41").Value = "0.999"

$ws.Range("D42").Value = "151.62"
$ws.Range("E42").Value = "  +0.16%  "

$ws.Range("D43").Value = "3.84"
$ws.Range("E43").Value = "  +2.86%  "

$ws.Range("D44").Value = "21.19"
$ws.Range("E44").Value = "  +3.12%  "

$ws.Range("E45").Value = "  +7.23%  "

$ws.Range("D46").Value = "0.617"
$ws.Range("E46").Value = "  +2.65%  "

$ws.Range("E47").Value = "  +2.92%  "

$ws.Range("D48").Value = "0.0241"
$ws.Range("E48").Value = "  +2.63%  "

$ws.Range("D49").Value = "18.53"
$ws.Range("E49").Value = "  +1.31%  "

$ws.Range("D50").Value = "1.75"
$ws.Range("E50").Value = "  -0.72%  "

$ws.Range("E51").Value = "  -0.23%  "
